# Corrects several "logic problem" values in column A (and their
# downstream column C totals) on sheet "données07", per the commit
# "modified data (because there were some logic problems)".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 8.2799999999999994
$ws.Range("C16").Value = 99
$ws.Range("A22").Value = 30.61
$ws.Range("C22").Value = 86
$ws.Range("A24").Value = 37.75
$ws.Range("C24").Value = 97
$ws.Range("A25").Value = 29.32
$ws.Range("C25").Value = 84
$ws.Range("A30").Value = 36.130000000000003
$ws.Range("C30").Value = 96
$ws.Range("A31").Value = 61.39
$ws.Range("C31").Value = 96
$ws.Range("A32").Value = 47.75
$ws.Range("C32").Value = 87
$ws.Range("A34").Value = 56.04
$ws.Range("C34").Value = 101
$ws.Range("A41").Value = 6.4399999999999995
$ws.Range("C41").Value = 79
$ws.Range("A44").Value = 6.2600000000000007
$ws.Range("C44").Value = 89
$ws.Range("A45").Value = 23.66
$ws.Range("C45").Value = 94
$ws.Range("A50").Value = 22.73
$ws.Range("C50").Value = 98
$ws.Range("A52").Value = 57.769999999999996
$ws.Range("C52").Value = 95
$ws.Range("A53").Value = 72.850000000000009
$ws.Range("C53").Value = 97
$ws.Range("A57").Value = 11.15
$ws.Range("C57").Value = 89
$ws.Range("A61").Value = 23.26
$ws.Range("C61").Value = 92
$ws.Range("A63").Value = 21.13
$ws.Range("C63").Value = 83
